# Insert two new rows of "Pimiento" data at row 681 in the Hortaliza /
# Macroferia Regional de Talca sheet, pushing the existing rows 681-774
# down to 683-776 (dimension grows from A1:R774 to A1:R776).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the block that starts at row 681.
$ws.Rows("681:682").Insert()

# ---- New row 681 ----
$ws.Cells.Item(681, 1).Value = 5
$ws.Cells.Item(681, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(681, 3).Value = "Maule"
$ws.Cells.Item(681, 4).Value = 44984
$ws.Cells.Item(681, 5).Value = 7
$ws.Cells.Item(681, 6).Value = 100112002
$ws.Cells.Item(681, 7).Value = "Pimiento"
$ws.Cells.Item(681, 8).Value = "Cuatro cascos rojo"
$ws.Cells.Item(681, 9).Value = "Primera"
$ws.Cells.Item(681, 10).Value = 150
$ws.Cells.Item(681, 11).Value = 13000
$ws.Cells.Item(681, 12).Value = 13000
$ws.Cells.Item(681, 13).Value = 13000
$ws.Cells.Item(681, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(681, 15).Value = "Región del Maule"
$ws.Cells.Item(681, 16).Value = 722
$ws.Cells.Item(681, 17).Value = 18
$ws.Cells.Item(681, 18).Value = "Hortaliza"

# ---- New row 682 ----
$ws.Cells.Item(682, 1).Value = 5
$ws.Cells.Item(682, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(682, 3).Value = "Maule"
$ws.Cells.Item(682, 4).Value = 44984
$ws.Cells.Item(682, 5).Value = 7
$ws.Cells.Item(682, 6).Value = 100112002
$ws.Cells.Item(682, 7).Value = "Pimiento"
$ws.Cells.Item(682, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(682, 9).Value = "Primera"
$ws.Cells.Item(682, 10).Value = 300
$ws.Cells.Item(682, 11).Value = 9000
$ws.Cells.Item(682, 12).Value = 9000
$ws.Cells.Item(682, 13).Value = 9000
$ws.Cells.Item(682, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(682, 15).Value = "Región del Maule"
$ws.Cells.Item(682, 16).Value = 500
$ws.Cells.Item(682, 17).Value = 18
$ws.Cells.Item(682, 18).Value = "Hortaliza"

# Make sure the date cells keep the workbook's date number format (style
# index 2 == "YYYY-MM-DD HH:MM:SS"), matching every other row in column D.
$ws.Cells.Item(681, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(682, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
